$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column N header from "Number of trials" to "NumTrialsInABlock"
$ws.Range("N1").Value = "NumTrialsInABlock"

# Add an (empty, but formatted like the header row) cell in P1 - copy
# the header style from O1 so it reuses the existing bold style (s=1)
# rather than minting a new one, then clear its content.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").ClearContents()

# Shrink column O (Comments) width - the NOTES block below the table is
# being removed, so the wide comments column is no longer needed.
$ws.Columns.Item(15).ColumnWidth = 34.0

# Remove the NOTES block (rows 18-21) text, leaving the (formatted) cells
# in column A blank.
$ws.Range("A18").ClearContents()
$ws.Range("A19").ClearContents()
$ws.Range("A20").ClearContents()
$ws.Range("A21").ClearContents()
